# Sample Project / Main.xlsx — rule R40's "return" cell (B11) is retyped
# from the text "R40" to the text "1". The cell keeps its existing
# formatting (style) and remains a plain text value (not a number), so we
# stage the new text in a scratch cell formatted as Text, copy it, and
# paste *values only* into B11 — this swaps the stored string without
# disturbing B11's number format / style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "1"
$scratch.Copy()

$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues

$scratch.Clear()
$excel.CutCopyMode = $false
